# Update feed result logic
# Applies edits to the "Test Result" sheet of the smoke test scenarios workbook:
#  - Flip the "Receive" (C8) and "Prepare" (C9) result flags from TRUE to FALSE
#  - Extend the bordered "Remarks" column formatting (D) down to the now
#    still-bordered rows 10 and 11
#  - Move the active selection to C11
#  - Let the summary formula in E1 recalculate against the new values

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Result")
$ws.Activate()

# Flip the two feed-result booleans that drive the summary counts
$ws.Range("C8").Value = $false
$ws.Range("C9").Value = $false

# Rows 10 and 11 need the same empty "Remarks" cell formatting (style) that
# the rows above them (D8/D9) already carry. Copy the formatting down.
$ws.Range("D9").Copy()
$ws.Range("D10").PasteSpecial(-4122)
$ws.Range("D11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Recalculate so the cached value of the E1 CONCATENATE/COUNTIF summary
# formula reflects the updated TRUE/FALSE counts
$excel.CalculateFull()

# Update the selected cell to match the saved view state
$ws.Range("C11").Select()
